$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly report rows (row 4 = week of 2022-11-09, row 5 = week of 2022-11-16)
# need to be swapped so the most recent week appears first (row 4) and the
# earlier week appears second (row 5). Swap the values for each affected
# column: D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).

$cols = @("D", "L", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $cell4 = $ws.Range("$col" + "4")
    $cell5 = $ws.Range("$col" + "5")

    $val4 = $cell4.Value()
    $val5 = $cell5.Value()

    $cell4.Value = $val5
    $cell5.Value = $val4
}
